$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = New-Object "object[,]" 1,20
$row[0,0] = "ECs"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "ECs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 7.567183666666668
$row[0,7] = 22.701551
$row[0,8] = 0.2104710071241515
$row[0,9] = 0.2104710071241515
$row[0,10] = 1
$row[0,11] = 0.3333333333333333
$row[0,12] = 0.5610396666666667
$row[0,13] = 1.683119
$row[0,14] = 0.04888859506813145
$row[0,15] = 0.04888859506813144
$row[0,16] = 4.245490201952112
$row[0,17] = 38.20941181756901
$row[0,18] = 0.01028963184087445
$row[0,19] = 0.01028963184087445
$ws.Range("A2:T2").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "ECs"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "FAPs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 7.567183666666668
$row[0,7] = 22.701551
$row[0,8] = 0.2104710071241515
$row[0,9] = 0.2104710071241515
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 7.236132333333333
$row[0,13] = 21.708397
$row[0,14] = 0.6305513932830891
$row[0,15] = 0.630551393283089
$row[0,16] = 54.75714240263856
$row[0,17] = 492.814281623747
$row[0,18] = 0.1327127867878287
$row[0,19] = 0.1327127867878287
$ws.Range("A3:T3").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "ECs"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "sCs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 7.567183666666668
$row[0,7] = 22.701551
$row[0,8] = 0.2104710071241515
$row[0,9] = 0.2104710071241515
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.678708333333333
$row[0,13] = 11.036125
$row[0,14] = 0.3205600116487796
$row[0,15] = 0.3205600116487795
$row[0,16] = 27.83746161443056
$row[0,17] = 250.537154529875
$row[0,18] = 0.06746858849544839
$row[0,19] = 0.06746858849544836
$ws.Range("A4:T4").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "FAPs"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "ECs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 8.450985000000001
$row[0,7] = 25.352955
$row[0,8] = 0.2350527491457871
$row[0,9] = 0.2350527491457871
$row[0,10] = 1
$row[0,11] = 0.3333333333333333
$row[0,12] = 0.5610396666666667
$row[0,13] = 1.683119
$row[0,14] = 0.04888859506813145
$row[0,15] = 0.04888859506813144
$row[0,16] = 4.741337807405001
$row[0,17] = 42.67204026664501
$row[0,18] = 0.01149139867263947
$row[0,19] = 0.01149139867263946
$ws.Range("A5:T5").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "FAPs"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "FAPs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 8.450985000000001
$row[0,7] = 25.352955
$row[0,8] = 0.2350527491457871
$row[0,9] = 0.2350527491457871
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 7.236132333333333
$row[0,13] = 21.708397
$row[0,14] = 0.6305513932830891
$row[0,15] = 0.630551393283089
$row[0,16] = 61.15244580701501
$row[0,17] = 550.3720122631349
$row[0,18] = 0.1482128384688965
$row[0,19] = 0.1482128384688965
$ws.Range("A6:T6").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "FAPs"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "sCs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 8.450985000000001
$row[0,7] = 25.352955
$row[0,8] = 0.2350527491457871
$row[0,9] = 0.2350527491457871
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.678708333333333
$row[0,13] = 11.036125
$row[0,14] = 0.3205600116487796
$row[0,15] = 0.3205600116487795
$row[0,16] = 31.088708944375
$row[0,17] = 279.798380499375
$row[0,18] = 0.07534851200425119
$row[0,19] = 0.07534851200425116
$ws.Range("A7:T7").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "M1"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "ECs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 9.443581333333334
$row[0,7] = 28.330744
$row[0,8] = 0.2626604773504909
$row[0,9] = 0.2626604773504908
$row[0,10] = 1
$row[0,11] = 0.3333333333333333
$row[0,12] = 0.5610396666666667
$row[0,13] = 1.683119
$row[0,14] = 0.04888859506813145
$row[0,15] = 0.04888859506813144
$row[0,16] = 5.29822372339289
$row[0,17] = 47.684013510536
$row[0,18] = 0.01284110171759026
$row[0,19] = 0.01284110171759026
$ws.Range("A8:T8").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "M1"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "FAPs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 9.443581333333334
$row[0,7] = 28.330744
$row[0,8] = 0.2626604773504909
$row[0,9] = 0.2626604773504908
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 7.236132333333333
$row[0,13] = 21.708397
$row[0,14] = 0.6305513932830891
$row[0,15] = 0.630551393283089
$row[0,16] = 68.33500422859645
$row[0,17] = 615.015038057368
$row[0,18] = 0.1656209299537533
$row[0,19] = 0.1656209299537532
$ws.Range("A9:T9").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "M1"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "sCs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 9.443581333333334
$row[0,7] = 28.330744
$row[0,8] = 0.2626604773504909
$row[0,9] = 0.2626604773504908
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.678708333333333
$row[0,13] = 11.036125
$row[0,14] = 0.3205600116487796
$row[0,15] = 0.3205600116487795
$row[0,16] = 34.74018134744445
$row[0,17] = 312.661632127
$row[0,18] = 0.08419844567914735
$row[0,19] = 0.08419844567914732
$ws.Range("A10:T10").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "M2"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "ECs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 6.387937333333333
$row[0,7] = 19.163812
$row[0,8] = 0.1776718609216568
$row[0,9] = 0.1776718609216568
$row[0,10] = 1
$row[0,11] = 0.3333333333333333
$row[0,12] = 0.5610396666666667
$row[0,13] = 1.683119
$row[0,14] = 0.04888859506813145
$row[0,15] = 0.04888859506813144
$row[0,16] = 3.583886232180889
$row[0,17] = 32.254976089628
$row[0,18] = 0.008686127663600249
$row[0,19] = 0.008686127663600249
$ws.Range("A11:T11").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "M2"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "FAPs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 6.387937333333333
$row[0,7] = 19.163812
$row[0,8] = 0.1776718609216568
$row[0,9] = 0.1776718609216568
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 7.236132333333333
$row[0,13] = 21.708397
$row[0,14] = 0.6305513932830891
$row[0,15] = 0.630551393283089
$row[0,16] = 46.22395988104044
$row[0,17] = 416.0156389293639
$row[0,18] = 0.1120312394513499
$row[0,19] = 0.1120312394513499
$ws.Range("A12:T12").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "M2"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "sCs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 6.387937333333333
$row[0,7] = 19.163812
$row[0,8] = 0.1776718609216568
$row[0,9] = 0.1776718609216568
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.678708333333333
$row[0,13] = 11.036125
$row[0,14] = 0.3205600116487796
$row[0,15] = 0.3205600116487795
$row[0,16] = 23.49935830094444
$row[0,17] = 211.4942247085
$row[0,18] = 0.05695449380670666
$row[0,19] = 0.05695449380670665
$ws.Range("A13:T13").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "sCs"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "ECs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 4.103880666666666
$row[0,7] = 12.311642
$row[0,8] = 0.1141439054579135
$row[0,9] = 0.1141439054579135
$row[0,10] = 1
$row[0,11] = 0.3333333333333333
$row[0,12] = 0.5610396666666667
$row[0,13] = 1.683119
$row[0,14] = 0.04888859506813145
$row[0,15] = 0.04888859506813144
$row[0,16] = 2.302439841266444
$row[0,17] = 20.721958571398
$row[0,18] = 0.005580335173427015
$row[0,19] = 0.005580335173427014
$ws.Range("A14:T14").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "sCs"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "FAPs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 4.103880666666666
$row[0,7] = 12.311642
$row[0,8] = 0.1141439054579135
$row[0,9] = 0.1141439054579135
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 7.236132333333333
$row[0,13] = 21.708397
$row[0,14] = 0.6305513932830891
$row[0,15] = 0.630551393283089
$row[0,16] = 29.69622358420822
$row[0,17] = 267.266012257874
$row[0,18] = 0.07197359862126058
$row[0,19] = 0.07197359862126056
$ws.Range("A15:T15").Value = $row

$row = New-Object "object[,]" 1,20
$row[0,0] = "sCs"
$row[0,1] = "Hras"
$row[0,2] = "Agtr1a"
$row[0,3] = "sCs"
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 4.103880666666666
$row[0,7] = 12.311642
$row[0,8] = 0.1141439054579135
$row[0,9] = 0.1141439054579135
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 3.678708333333333
$row[0,13] = 11.036125
$row[0,14] = 0.3205600116487796
$row[0,15] = 0.3205600116487795
$row[0,16] = 15.09698000747222
$row[0,17] = 135.87282006725
$row[0,18] = 0.03658997166322596
$row[0,19] = 0.03658997166322595
$ws.Range("A16:T16").Value = $row
